# Fixed duplicate "2019" data values
# The weekly-incidence dates in column B of the "Incidence" sheet had a
# block of rows mistakenly stamped with 2019 dates that duplicated later
# rows. Correct them back to the matching 2018 dates (same month/day,
# 365 days earlier) so the series is a continuous weekly run into 2019.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Incidence")

$corrections = @{
    2  = 43317
    3  = 43324
    4  = 43332
    5  = 43338
    6  = 43345
    8  = 43359
    9  = 43366
    10 = 43375
    11 = 43380
    12 = 43388
    13 = 43394
    14 = 43401
    15 = 43408
    16 = 43415
    17 = 43423
    18 = 43430
    19 = 43437
    20 = 43444
    21 = 43450
    22 = 43459
}

foreach ($row in $corrections.Keys) {
    $ws.Cells.Item($row, 2).Value = $corrections[$row]
}

# Move the active selection on the frozen-pane view to B22, matching the
# last-edited cell.
$ws.Activate()
$ws.Range("B22").Select()
